# Auto-generated from XML diff: updates Price/Volume(1h) figures (and a
# KickToken/CEJI row swap) in the crypto symbol-list sheet, matching the
# "Updated symbol list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold text-formatted values (e.g.
# "245.99", "0.01%") in the source workbook, not numbers. Force the number
# format to Text first so Excel doesn't quietly convert them to numeric/
# percentage cells, which would change the stored cell type.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "246.02"
$ws.Range("E2").Value = "0.16%"
$ws.Range("D3").Value = "28.41"
$ws.Range("E3").Value = "-2.42%"
$ws.Range("D4").Value = "5.291"
$ws.Range("E4").Value = "2.15%"
$ws.Range("E5").Value = "-0.39%"
$ws.Range("D6").Value = "6.645"
$ws.Range("E6").Value = "1.24%"
$ws.Range("D7").Value = "3.213"
$ws.Range("E7").Value = "3.36%"
$ws.Range("D8").Value = "0.8625"
$ws.Range("E8").Value = "0.22%"
$ws.Range("D9").Value = "0.8855"
$ws.Range("E9").Value = "2.80%"
$ws.Range("D10").Value = "0.1386"
$ws.Range("E10").Value = "1.57%"
$ws.Range("D11").Value = "0.07087"
$ws.Range("E11").Value = "0.31%"
$ws.Range("D12").Value = "0.03149"
$ws.Range("D13").Value = "0.09231"
$ws.Range("E13").Value = "-1.55%"
$ws.Range("E14").Value = "-1.02%"
$ws.Range("D15").Value = "0.0005967"
$ws.Range("E15").Value = "-94.18%"
$ws.Range("D16").Value = "0.006019"
$ws.Range("E16").Value = "0.33%"
$ws.Range("D17").Value = "3.497"
$ws.Range("E17").Value = "0.20%"
$ws.Range("E18").Value = "-4.60%"
$ws.Range("D19").Value = "0.3122"
$ws.Range("E19").Value = "-2.46%"
$ws.Range("E20").Value = "0.60%"
$ws.Range("E21").Value = "2.08%"
$ws.Range("D22").Value = "3.489"
$ws.Range("E22").Value = "0.42%"
$ws.Range("D23").Value = "0.04097"
$ws.Range("E23").Value = "-1.37%"
$ws.Range("D24").Value = "0.1378"
$ws.Range("E24").Value = "-0.16%"
$ws.Range("D25").Value = "0.001219"
$ws.Range("E25").Value = "-0.67%"
$ws.Range("D26").Value = "0.004167"
$ws.Range("E26").Value = "-16.47%"
$ws.Range("E27").Value = "-0.97%"
$ws.Range("D40").Value = "0.03813"
$ws.Range("E40").Value = "1.98%"
$ws.Range("D41").Value = "0.1069"
$ws.Range("E41").Value = "0.02%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002199"
$ws.Range("E42").Value = "4.61%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.002948"
$ws.Range("E43").Value = "-49.03%"
$ws.Range("D44").Value = "0.009476"
$ws.Range("E44").Value = "10.02%"
$ws.Range("D45").Value = "0.00005279"
$ws.Range("E45").Value = "-0.05%"
$ws.Range("E46").Value = "-0.14%"
$ws.Range("D47").Value = "0.08906"
$ws.Range("E47").Value = "56.10%"
$ws.Range("D48").Value = "0.002259"
$ws.Range("E48").Value = "0.08%"
$ws.Range("E49").Value = "-0.14%"
$ws.Range("E50").Value = "-0.14%"
